$d = $word.ActiveDocument

# --- Change 1: "Sprint Number: 3" -> "Sprint Number: 2" ---
$find1 = $d.Content
$find1.Find.Execute("Sprint Number: 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$numRange = $d.Range($find1.End - 1, $find1.End)
$numRange.Text = "2"

# --- Change 2: "Writing unit tests" -> "Analyze the functions to write unit tests" ---
$find2 = $d.Content
$find2.Find.Execute("Writing unit tests", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find2.Text = "Analyze the functions to write unit tests"

# --- Change 3: "Continue writing unit tests" -> "Analyze the functions to write unit tests" ---
$find3 = $d.Content
$find3.Find.Execute("Continue writing unit tests", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find3.Text = "Analyze the functions to write unit tests"
